# Fill in the missing benchmark rows (8-13) and the CPU-percentage column (F)
# on the "Carga 400 - 1 thread" sheet, then apply Percentage formatting to F4:F13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Carga 400 - 1 thread")

# C (tiempo autenticacion), D (tiempo actualizacion), E (#transacciones perdidas)
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0

$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 1

$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0

$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 1

$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0

$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# F (Porcentaje CPU)
$ws.Range("F4").Value = 0.1095
$ws.Range("F5").Value = 0.112
$ws.Range("F6").Value = 0.1045
$ws.Range("F7").Value = 0.1033
$ws.Range("F8").Value = 0.1089
$ws.Range("F9").Value = 0.1077
$ws.Range("F10").Value = 0.0987
$ws.Range("F11").Value = 0.1012
$ws.Range("F12").Value = 0.1019
$ws.Range("F13").Value = 0.1021

$ws.Range("F4:F13").NumberFormat = "0.00%"
$ws.Range("F4:F13").Style = "Percent"

$ws.Select()
$ws.Range("I13").Select()
